$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the "Apparent Difficulty" ratings for the 'Going to Town' (row 10)
# and 'Mambo Marie' (row 11) tasks, and update 'Guns' (row 13) to match,
# as part of further explaining the 'Going to Town' design rationale.
$ws.Range("C10").Value = "Medium"
$ws.Range("C11").Value = "Hard"
$ws.Range("C13").Value = "Hard"

# Move the active selection to C12
$ws.Range("C12").Select()
